$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values that changed between the old and new layout
$ws.Range("B10").Value = "6495737 - Durval Rodrigues Junior"
$ws.Range("C10").Value = "6495737 - Durval Rodrigues Junior"
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "01/01/2012"
$ws.Range("C15").Value = "01/01/2012"
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "6495737 - Durval Rodrigues Junior"
$ws.Range("C18").Value = "6495737 - Durval Rodrigues Junior"
$ws.Range("A19").Value = "Critério:"
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B23").Value = "LOB1021 -  Física IV  (Requisito)`n"
$ws.Range("C23").Value = "LOB1021 -  Física IV  (Requisito)`n"
$ws.Range("B24").Value = "LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)`n"
$ws.Range("C24").Value = "LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)`n"
$ws.Range("B25").Value = "LOM3229 -  Métodos Experimentais da Física II  (Indicação de Conjunto)`n"
$ws.Range("C25").Value = "LOM3229 -  Métodos Experimentais da Física II  (Indicação de Conjunto)`n"

# Clear cells whose content was removed in the new layout
$ws.Range("B14").Clear()
$ws.Range("C14").Clear()
$ws.Range("B16").Clear()
$ws.Range("C16").Clear()
$ws.Range("B22").Clear()
$ws.Range("C22").Clear()
$ws.Range("A23").Clear()

# Adjust row heights to match the new layout
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
$ws.Rows.Item(23).RowHeight = 30
$ws.Rows.Item(17).AutoFit()
$ws.Rows.Item(22).AutoFit()

# Remove the now-unused trailing row
$ws.Rows.Item(26).Delete()

